# Update cryptos list data (price / 1h volume columns) to match the
# latest scrape, per commit "Updated cryptos list on Thu Aug 29
# 10:34:45 UTC 2024 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.420.94"
$ws.Range("E2").Value = "  -0.82%  "
$ws.Range("D3").Value = "2.520.04"
$ws.Range("E3").Value = "  -0.13%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "542.38"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.50"
$ws.Range("E6").Value = "  -1.14%  "
$ws.Range("E7").Value = "  -0.47%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.574"
$ws.Range("E8").Value = "  -0.76%  "
$ws.Range("D9").Value = "2.541.70"
$ws.Range("E9").Value = "  +0.69%  "
$ws.Range("E10").Value = "  -0.37%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.160"
$ws.Range("E11").Value = "  +0.27%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.59"
$ws.Range("E12").Value = "  +1.50%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.362"
$ws.Range("E13").Value = "  +1.87%  "
$ws.Range("D14").Value = "2.964.21"
$ws.Range("E14").Value = "  -0.11%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.57"
$ws.Range("E15").Value = "  -4.11%  "
$ws.Range("D16").Value = "59.331.96"
$ws.Range("E16").Value = "  -1.05%  "
$ws.Range("E17").Value = "  +1.00%  "
$ws.Range("D18").Value = "2.533.44"
$ws.Range("E18").Value = "  +1.21%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.19"
$ws.Range("E19").Value = "  -1.99%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.28"
$ws.Range("E20").Value = "  -2.12%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "326.50"
$ws.Range("E21").Value = "  -0.22%  "
$ws.Range("E22").Value = "  +0.27%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.91"
$ws.Range("E23").Value = "  +1.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "62.18"
$ws.Range("E24").Value = "  +0.94%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.433"
$ws.Range("E25").Value = "  -3.84%  "
$ws.Range("E26").Value = "  +1.63%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.995"
$ws.Range("E27").Value = "  -1.26%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.04"
$ws.Range("E28").Value = "  +2.89%  "
$ws.Range("E29").Value = "  -1.38%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.82"
$ws.Range("E30").Value = "  -0.24%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.72"
$ws.Range("E31").Value = "  -1.29%  "
$ws.Range("E32").Value = "  -6.21%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.48"
$ws.Range("E33").Value = "  +2.77%  "
$ws.Range("E34").Value = "  -0.14%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "158.86"
$ws.Range("E35").Value = "  +0.47%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.77"
$ws.Range("E36").Value = "  -1.25%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.38"
$ws.Range("E37").Value = "  -2.94%  "
$ws.Range("E38").Value = "  -7.08%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "37.03"
$ws.Range("E39").Value = "  +0.56%  "
$ws.Range("E40").Value = "  -6.09%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.835"
$ws.Range("E41").Value = "  -0.25%  "
$ws.Range("B42").Value = "Bittensor"
$ws.Range("C42").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "296.59"
$ws.Range("E42").Value = "  -6.30%  "
$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.69"
$ws.Range("E43").Value = "  -2.65%  "
$ws.Range("E44").Value = "  -0.63%  "
$ws.Range("E45").Value = "  -0.32%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.81"
$ws.Range("E46").Value = "  +0.87%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0937"
$ws.Range("E47").Value = "  -0.42%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.80"
$ws.Range("E48").Value = "  +0.34%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "123.27"
$ws.Range("E49").Value = "  -2.77%  "
$ws.Range("E50").Value = "  -2.29%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0514"
$ws.Range("E51").Value = "  -3.95%  "
